# Updates to Embodied Carbon sheet
# - Fill in manufacturing/process data for the SoC rows on the
#   "Embodied Carbon" sheet (CPU/GPU rows now say "USA" for
#   manufacturing location with lithography/EPA/GPA figures, and the
#   DRAM row gets a die area + "TSMC" manufacturing location with its
#   own CI_fab/lithography figures).
# - Leave a trail of view/selection state matching where the author was
#   last working (MIT Nodes G2, Embodied Carbon C11).

$wb = $excel.ActiveWorkbook

$wsNodes = $wb.Worksheets.Item("MIT Nodes")
$wsNodes.Activate()
$wsNodes.Range("G2").Select()

$ws = $wb.Worksheets.Item("Embodied Carbon")
$ws.Activate()

# Row 3 - CPU (Intel Xeon Gold 6248)
$ws.Range("D3").Value = "USA"
$ws.Range("E3").Value = 380
$ws.Range("F3").Value = 14
$ws.Range("G3").Value = 1.2
$ws.Range("I3").Value = 200

# Row 4 - GPU (Intel Xeon Platinum 8260)
$ws.Range("D4").Value = "USA"
$ws.Range("E4").Value = 380
$ws.Range("F4").Value = 14
$ws.Range("G4").Value = 1.2
$ws.Range("I4").Value = 200

# Row 5 - DRAM (Nvidia Volta V100)
$ws.Range("C5").Value = 8.15
$ws.Range("D5").Value = "TSMC"
$ws.Range("E5").Value = 583
$ws.Range("F5").Value = 12

$ws.Range("C11").Select()
